$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @{
    2  = @(0.05688009254515978, 0.3078030679106084, 0.1763561601164249, 0.4199478064193512, 0.420389698777502, 49)
    3  = @(0.2667693795774359, 0.5466263928059274, 0.7076554788746405, 0.8412226095835992, 0.8060705778746855, 49)
    4  = @(0.5431296193130369, 0.7231827691434326, 1.246365227305445, 1.116407285584184, 0.985706586491255, 48)
    5  = @(0.4692704090458715, 0.7347082924362607, 1.284144682252101, 1.133201077590425, 1.042621157342553, 47)
    6  = @(0.3883654118239238, 0.6070310945326703, 0.8839353291870329, 0.9401783496693767, 0.8656782987275211, 46)
    7  = @(0.3703187222185918, 0.5607740953733055, 0.5932288700479009, 0.7702135223740888, 0.685502670612526, 34)
    8  = @(0.3586222435103067, 0.5485688517592353, 0.6215876060107235, 0.788408273682312, 0.7130099896635785, 33)
    9  = @(0.1822978435501495, 0.5070926453844893, 0.3970354312988059, 0.6301074759902519, 0.6229417758010422, 16)
    10 = @(0.203264939200266, 0.4823557971524242, 0.379312118145266, 0.6158832017073254, 0.6128217818122421, 10)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
